# Fill in the "Column2" calculated column of Table1 (Sheet1) with the
# formula:   =Table1[[#This Row],[Pos*5]]-E3
# Setting the Formula across the whole DataBodyRange at once lets Excel
# auto-adjust the relative "E3" reference per row (E3, E4, E5, ... down to
# E43), exactly as if the formula had been typed into the first data row of
# the table column and auto-filled down the rest of the calculated column.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

$col2 = $tbl.ListColumns.Item("Column2")
[void]($col2.DataBodyRange.Formula = "=Table1[[#This Row],[Pos*5]]-E3")

# Restore the recorded selection/scroll state: active cell J19, sheet
# scrolled back so row 1 is the top-left visible cell again.
[void]$ws.Activate()
[void]$ws.Range("J19").Select()
